# D3_study_population_target_cohorts.xlsx -- "study population cohort addition"
#
# The "Data Model" sheet documents the cohort_* variables (is_in_*,
# cohort_entry_date_*, cohort_exit_date_*, is_censored_in_*) in rows 8-11.
# Their free-text notes in column K were previously shifted down by one row
# (row 9's note belonged to row 8's variable, etc.) and the wording of two
# notes was tightened up ("september (???) XXXX" -> "1st september XXXX",
# "(april?)" -> "30th april XXXX"). This realigns the notes with the correct
# row and removes the stray leftover note in row 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Model")

# --- Column K notes, realigned to the row describing the relevant variable ---

# Row 10: cohort_exit_date_cohort_type_and_label -- gets the "what date is
# the exit date" note (previously mis-placed on row 11), with the
# seasonalXXX wording tightened up.
$ws.Range("K10").Value = "for birth cohorts: for birthYY, exit date is earliest between the date they turn XX months old and study_exit_date`nfor adolescence: exit date is earliest between the date they turn 16 years old and study_exit_date`nfor seasonalXXX: exit date is earliest between end of season 30th april XXXX and study_exit_date`nfor covid_vacc: study_exit_date"

# Row 8: is_in_cohort_type_and_label -- gets the "when is a person in the
# cohort" note (previously mis-placed on row 9).
$ws.Range("K8").Value = "for birth cohorts: the person is in the cohort if birth_date is between study_entry_date and study_exit_date`nfor adolescence: the person is in the cohort if the ninth birthday is between study_entry_date and study_exit_date`nfor seasonalXXX: the person is in the cohort if 1st september XXXX is between study_entry_date and study_exit_date`nforcovid_vacc: the person is in the cohort if 1st december 2020 is between study_entry_date and study_exit_date`n"

# Row 9: cohort_entry_date_cohort_type_and_label -- gets the "what date is
# the entry date" note (previously mis-placed on row 10).
$ws.Range("K9").Value = "for birth cohorts: birth_date`nfor adolescence: ninth birthday `nfor seasonalXXX: september (???) XXXX `nfor covid_vacc: 1st december 2020"

# Row 11: is_censored_in_cohort_type_and_label -- no note of its own anymore;
# clear the stray leftover text that used to sit here.
$ws.Range("K11").ClearContents()

# --- Row heights follow the (now longer/shorter) wrapped note text ---
$ws.Rows.Item(8).RowHeight = 225
$ws.Rows.Item(9).RowHeight = 62.25
$ws.Rows.Item(10).RowHeight = 150
$ws.Rows.Item(11).AutoFit()

# --- View state: the author ended up on the "Data Model" sheet, zoomed in a
# touch more, and left the selection on A11; the "Example" sheet's
# selection moved to L1. ---
$exws = $wb.Worksheets.Item("Example")
$exws.Range("L1").Select()

$ws.Activate()
$ws.Range("A11").Select()
$excel.ActiveWindow.Zoom = 130

# --- Cosmetic locale fix: the base cell style was saved under its Italian
# name ("Normale"); rename it back to the English default.
$style = $wb.Styles.Item("Normale")
$style.Name = "Normal"
